$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509960664747276"
$ws1.Range("B2").Value = "go_stims-16509960664426973.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960664586995.csv"
$ws1.Range("B4").Value = "go_stims-16509960664586995.csv"
$ws1.Range("B5").Value = "GNG_stims-16509960664747276.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509960690186963"
$ws2.Range("B2").Value = "ZB-match_9-16509960665946972.csv"
$ws2.Range("B3").Value = "OB-16509960677227042.csv"
$ws2.Range("B4").Value = "TB-16509960689947357.csv"
$ws2.Range("B5").Value = "TB-1650996068770696.csv"
$ws2.Range("B6").Value = "OB-16509960680267055.csv"
$ws2.Range("B7").Value = "OB-16509960670587342.csv"
$ws2.Range("B8").Value = "ZB-match_3-16509960666427047.csv"
$ws2.Range("B9").Value = "ZB-match_5-16509960665227375.csv"
$ws2.Range("B10").Value = "TB-1650996068818698.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509960690186963"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509960690747051"
$ws4.Range("B2").Value = "MM_stims-16509960690427377.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960690186963.csv"
$ws4.Range("B4").Value = "MM_stims-16509960690587347.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960690427377.csv"
$ws4.Range("B6").Value = "MM_stims-16509960690747051.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960690587347.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650996069146698"
$ws5.Range("B2").Value = "vSAT_stims-16509960691147354.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960690987008.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509960691307352.csv"
$ws5.Range("B5").Value = "SAT_stims-16509960690747051.csv"
